$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update intro/instruction texts (shared strings referenced by D6 and D11) ---
$ws.Range("D6").Value = '<p>Today we are going to learn about <a href="/information/crafting" target="_blank">crafting</a>. There are various types of crafting, but we want to focus on weapons and Armour for now.</p><p>At first you will fail a lot, even with the suggested <a href="/information/quest-items" target="_blank">quest items</a> this quest will take roughly <strong>4 hours to complete</strong>. Do not worry about trying to rush through it. Instead focus on your gear that you get from exploration and if any of it increases your Weapon Crafting or Armour Crafting – equip it.</p><p>Eventually you will craft beyond what the <a href="/information/shop" target="_blank">shop</a> sells. The Shop stops at two billion gold, where as players can craft up to 36+ Billion gold items that are much better then shop gear.</p><p><strong>Crafting cannot be automated. Keep an eye on Server Messages section to see successes, failures and if you have new items to craft.</strong></p><p><strong>Desktop:</strong></p><p>- First we want to find the items we need, quest items are automatically used. You can open the Teleport Map Action to then select the location: Ruined Port City Of Kalize (X/Y): 32/368 and click teleport if you can afford the cost. This will get you the: Weapon Smith’s Book which adds 25% to Skill Bonus and XP.</p><p>- Repeat the above step to then go to: Dragon cliffs (X/Y): 192/176 to get the Blacksmith’s book for the same bonuses towards Armour Crafting.</p><p>Both of these can be upgraded later on when we have access to Labyrinth. There are some One Off <a href="/information/quests" target="_blank">quests</a> that will upgrade these.</p><p>- Now lets craft, first set up exploration for an hour or two – or what ever you feel like doing, Exploration will run while logged out. Remember to set it up with a monster you can kill in one hit.</p><p>- Close exploration and then from the drop down Craft/Enchant select craft, select weapons, select Broken Dagger and then click craft.</p><p><strong>Mobile:</strong></p><p>- Select Map from the actions drop down.</p><p>- Click Teleport from the actions under the map.</p><p>- Select the location: Ruined Port City Of Kalize (X/Y): 32/368 and click teleport if you can afford the cost. This will get you the: Weapon Smith’s Book which adds 25% to Skill Bonus and XP.</p><p>- Repeat the above step to then go to: Dragon cliffs (X/Y): 192/176 to get the Blacksmith’s book for the same bonuses towards Armour Crafting.</p><p>Both of these can be upgraded later on when we have access to Labyrinth. There are some One off <a href="/information/quests" target="_blank">quests</a> that will upgrade these.</p><p>- Now lets craft, first set up exploration for an hour or two – or what ever you feel like doing, Exploration will run while logged out. Remember to set it up with a monster you can kill in one hit.</p><p>- Close exploration, select Craft from the drop down and then select Craft from the Craft/Enchant</p><p>- Select Weapons, Broken Dagger – Click craft.</p><p>- When ready, click Change Type, select Armour, select an Armour to craft and click craft.</p><p>Over time new items will be added to the list. That’s all there is to it. When you are ready, click Change Type, select Armour, pick an item to craft and repeat.</p>'
$ws.Range("D11").Value = '<p>Now we learn about a new feature: <a href="/information/class-skills" target="_blank">Class Skills</a> and Class Bonus.</p><p>Class Skills are different for each <a href="/information/races-and-classes#3" target="_blank">class</a> in the game. Every class has a skill which you can see on your skills table under Training Tab. It has an icon beside it and is in orange text.</p><p>This skill important to level because it allows you to increase your Class Bonus which can be seen on the character sheet to the left, under inventory count or on mobile under Class Details, at the bottom of the details section.</p><p>Every class has a special attack that fires automatically based on Three Things:</p><p>- Class Bonus % (the higher, the more chance for the special to fire)</p><p>- Weapons</p><p>- Attack type</p><p>For example, Heretics:</p><p>With a damage spell equipped you have a small chance to cast another spell. Enemies cannot avoid this.</p><p>This means while casting and with at least one damage spell, based on the % of your class bonus you can cast another spell, automatically. Each class has its own special which you can read about in the help docs for your specific class.</p><p>To do this:</p><p><strong>Desktop/Mobile</strong></p><p>- The instructions state to level a Effects Class to the specified level. To do this, go to your character sheet section, in your skill section – train the skill with the orange text to the specified level.</p><p>- To Get the Gold Dust you can disenchant items that drop you no longer need, or craft some items, enchant them and disenchant them. Now is the time to explore the <a href="/information/enchanting" target="_blank">Enchanting</a> list in the docs to see what types of enchantments you can apply to your equipment so you can start creating a gear set geared towards your needs.</p><p>That’s it. Now you might not see your class bonus fire off much at first, but over time, keep leveling this skill and you will start to see your special fire off automatically when manually fighting.</p>'

# --- Row 11 edits ---
# Remove required_skill (F11) and required_skill_level (G11)
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()
# required_faction_id: Labyrinth -> Dungeons
$ws.Range("P11").Value = "Dungeons"
# required_stats: 500 -> 750
$ws.Range("AD11").Value = 750

# --- Row 12 edits ---
# required_skill: Casting Accuracy -> Weapon Crafting
$ws.Range("F12").Value = "Weapon Crafting"
# required_skill_level: 30 -> 25
$ws.Range("G12").Value = 25
# required_secondary_skill_level: 50 -> 25
$ws.Range("I12").Value = 25
# required_skill_type_level: 30 -> 50
$ws.Range("K12").Value = 50
# required_stats: 600 -> 900
$ws.Range("AD12").Value = 900

# --- Column D width grew because its text content grew (bestFit column) ---
$ws.Columns.Item(4).ColumnWidth = 3871.0
